$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.761.01"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.865.72"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").Value = "'1.038"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.05%  "
$ws.Range("D5").Value = "'323.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'1.033"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "'0.4428"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").Value = "'0.3803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "'0.07479"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").Value = "'0.8877"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("E12").Value = "  -6.12%  "
$ws.Range("D13").Value = "'5.553"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "'6.782"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'84.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "'0.000009144"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "'1.033"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "27.753.39"
$ws.Range("D22").Value = "'5.322"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").Value = "'11.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").Value = "2.092.21"
$ws.Range("E24").Value = "  -4.51%  "
$ws.Range("D25").Value = "'2.023"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.26%  "
$ws.Range("D26").Value = "'158.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'18.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "'5.354"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "'1.990"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("D30").Value = "'119.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("D31").Value = "'0.09070"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "'0.7794"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "'3.034"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'4.612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").Value = "'1.035"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "'1.147"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'0.01990"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.05365"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "'2.898"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").Value = "'0.5220"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "'0.1694"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").Value = "'6.906"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.72%  "
$ws.Range("D44").Value = "'8.717"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'110.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").Value = "'0.06731"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.00%  "
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'1.718"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("D50").Value = "'0.4732"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").Value = "'1.918"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.52%  "
